$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 15703.4
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 15703.4
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H51").Value = 2832.5715
$ws.Range("I51").Value = 2232.25
$ws.Range("J51").Value = 3633
$ws.Range("K51").Value = 2232.25
$ws.Range("L51").Value = 3633
$ws.Range("M51").Value = -1748.25
$ws.Range("N51").Value = -4601
$ws.Range("H62").Value = 8786.75
$ws.Range("I62").Value = 5756.2856
$ws.Range("K62").Value = 5756.2856
$ws.Range("M62").Value = -5132.2856
$ws.Range("H65").Value = 8786.75
$ws.Range("I65").Value = 5756.2856
$ws.Range("K65").Value = 28781.428
$ws.Range("M65").Value = -25661.428
$ws.Range("H120").Value = 35000
$ws.Range("J120").Value = 35000
$ws.Range("L120").Value = 35000
$ws.Range("N120").Value = -44676
$ws.Range("H137").Value = 1105.2787
$ws.Range("I137").Value = 633.6129
$ws.Range("K137").Value = 1900.8387
$ws.Range("M137").Value = 649.1613000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6403.905
$ws.Range("I2").Value = 7630.1177
$ws.Range("J2").Value = 1192.5
$ws.Range("K2").Value = 7630.1177
$ws.Range("L2").Value = 1192.5
$ws.Range("M2").Value = -7517.1177
$ws.Range("N2").Value = -1418.5
$ws.Range("H32").Value = 3392.6567
$ws.Range("I32").Value = 2198.0408
$ws.Range("J32").Value = 6644.6665
$ws.Range("K32").Value = 2198.0408
$ws.Range("L32").Value = 6644.6665
$ws.Range("M32").Value = -1911.0408
$ws.Range("N32").Value = -7218.6665
$ws.Range("H45").Value = 1777.1765
$ws.Range("I45").Value = 1393.7142
$ws.Range("J45").Value = 3566.6667
$ws.Range("K45").Value = 1393.7142
$ws.Range("L45").Value = 3566.6667
$ws.Range("M45").Value = -1016.7142
$ws.Range("N45").Value = -4320.6667
$ws.Range("H61").Value = 1528.1515
$ws.Range("I61").Value = 1209
$ws.Range("J61").Value = 2525.5
$ws.Range("K61").Value = 1209
$ws.Range("L61").Value = 2525.5
$ws.Range("M61").Value = -997
$ws.Range("N61").Value = -2949.5
$ws.Range("H116").Value = 6403.905
$ws.Range("I116").Value = 7630.1177
$ws.Range("J116").Value = 1192.5
$ws.Range("K116").Value = 7630.1177
$ws.Range("L116").Value = 1192.5
$ws.Range("M116").Value = -5336.1177
$ws.Range("N116").Value = -5780.5
$ws.Range("H132").Value = 2161.1428
$ws.Range("I132").Value = 1798.1034
$ws.Range("J132").Value = 3915.8333
$ws.Range("K132").Value = 5394.3102
$ws.Range("L132").Value = 11747.4999
$ws.Range("M132").Value = -2864.3102
$ws.Range("N132").Value = -16807.4999
$ws.Range("H136").Value = 1528.1515
$ws.Range("I136").Value = 1209
$ws.Range("J136").Value = 2525.5
$ws.Range("K136").Value = 3627
$ws.Range("L136").Value = 7576.5
$ws.Range("M136").Value = -1077
$ws.Range("N136").Value = -12676.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6403.905
$ws.Range("I3").Value = 7630.1177
$ws.Range("J3").Value = 1192.5
$ws.Range("K3").Value = 7630.1177
$ws.Range("L3").Value = 1192.5
$ws.Range("M3").Value = -7516.1177
$ws.Range("N3").Value = -1420.5
$ws.Range("H94").Value = 577.3333
$ws.Range("I94").Value = 599.8484999999999
$ws.Range("J94").Value = 329.66666
$ws.Range("K94").Value = 599.8484999999999
$ws.Range("L94").Value = 329.66666
$ws.Range("M94").Value = -148.8484999999999
$ws.Range("N94").Value = -1231.66666
$ws.Range("H134").Value = 3040.8857
$ws.Range("I134").Value = 1392.5
$ws.Range("J134").Value = 5830.4614
$ws.Range("K134").Value = 4177.5
$ws.Range("L134").Value = 17491.3842
$ws.Range("M134").Value = -1642.5
$ws.Range("N134").Value = -22561.3842

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1388.1
$ws.Range("I31").Value = 932.119
$ws.Range("K31").Value = 932.119
$ws.Range("M31").Value = -637.119
$ws.Range("H34").Value = 1388.1
$ws.Range("I34").Value = 932.119
$ws.Range("K34").Value = 932.119
$ws.Range("M34").Value = -730.119
$ws.Range("H94").Value = 515.2069
$ws.Range("I94").Value = 700.5
$ws.Range("J94").Value = 485.56
$ws.Range("K94").Value = 700.5
$ws.Range("L94").Value = 485.56
$ws.Range("M94").Value = -249.5
$ws.Range("N94").Value = -1387.56

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 853.64105
$ws.Range("I107").Value = 165.88889
$ws.Range("J107").Value = 1443.1428
$ws.Range("K107").Value = 497.66667
$ws.Range("L107").Value = 4329.428400000001
$ws.Range("M107").Value = 1422.33333
$ws.Range("N107").Value = -8169.428400000001
$ws.Range("H122").Value = 710.3333
$ws.Range("J122").Value = 1208.8572
$ws.Range("L122").Value = 10879.7148
$ws.Range("N122").Value = -15779.7148

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1164.9348
$ws.Range("I102").Value = 1045.2162
$ws.Range("J102").Value = 1657.1111
$ws.Range("K102").Value = 1045.2162
$ws.Range("L102").Value = 1657.1111
$ws.Range("M102").Value = 576.7837999999999
$ws.Range("N102").Value = -4901.1111
$ws.Range("H132").Value = 2201.1458
$ws.Range("I132").Value = 1787.4166
$ws.Range("J132").Value = 3442.3333
$ws.Range("K132").Value = 5362.2498
$ws.Range("L132").Value = 10326.9999
$ws.Range("M132").Value = -2832.2498
$ws.Range("N132").Value = -15386.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3810.15
$ws.Range("I122").Value = 3640.6
$ws.Range("J122").Value = 3866.6667
$ws.Range("K122").Value = 10921.8
$ws.Range("L122").Value = 11600.0001
$ws.Range("M122").Value = -8471.799999999999
$ws.Range("N122").Value = -16500.0001
$ws.Range("H132").Value = 3348.0435
$ws.Range("I132").Value = 2433.5312
$ws.Range("J132").Value = 5438.357
$ws.Range("K132").Value = 7300.5936
$ws.Range("L132").Value = 16315.071
$ws.Range("M132").Value = -4770.5936
$ws.Range("N132").Value = -21375.071

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 5051517
$ws.Range("I107").Value = 11112176
$ws.Range("J107").Value = 968.5
$ws.Range("K107").Value = 33336528
$ws.Range("L107").Value = 2905.5
$ws.Range("M107").Value = -33334608
$ws.Range("N107").Value = -6745.5
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

